# Apply updated quarterly margin figures to the INVH income sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INVH")

# Row 13 - Gross Margin
$ws.Range("B13").Value = 0.6349
$ws.Range("D13").Value = 0.6155
$ws.Range("E13").Value = 0.5998
$ws.Range("F13").Value = 0.5881
$ws.Range("G13").Value = 0.5854

# Row 14 - EBIT Margin
$ws.Range("B14").Value = 0.2725
$ws.Range("D14").Value = 0.2621
$ws.Range("E14").Value = 0.2527
$ws.Range("F14").Value = 0.2452
$ws.Range("G14").Value = 0.2303

# Row 15 - EBT margin
$ws.Range("D15").Value = 0.061
$ws.Range("F15").Value = 0.0464

# Row 16 - Net Profit Margin
$ws.Range("D16").Value = 0.0784
$ws.Range("E16").Value = 0.0866
$ws.Range("F16").Value = 0.0898
$ws.Range("G16").Value = 0.0822

# Row 17 - Free Cash Flow Margin
$ws.Range("B17").Value = 0.4956
$ws.Range("D17").Value = 0.653
$ws.Range("E17").Value = 0.6576
$ws.Range("F17").Value = 0.694
$ws.Range("G17").Value = 0.7344

# Row 28 - EBITDA Margin
$ws.Range("B28").Value = 0.6172
$ws.Range("D28").Value = 0.6026
$ws.Range("E28").Value = 0.5892
$ws.Range("F28").Value = 0.578
$ws.Range("G28").Value = 0.5651

# Row 29 - Operating Cash Flow Margin
$ws.Range("B29").Value = 0.436
$ws.Range("D29").Value = 0.4027
$ws.Range("E29").Value = 0.3897
$ws.Range("F29").Value = 0.3651
$ws.Range("G29").Value = 0.3752
